$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.046.10'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.871.27'
$ws.Range("E3").Value = '  -2.06%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.52'
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5036'
$ws.Range("E7").Value = '  -1.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3842'
$ws.Range("E8").Value = '  -2.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08647'
$ws.Range("E9").Value = '  -6.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.114'
$ws.Range("E10").Value = '  -2.23%  '
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.305'
$ws.Range("E12").Value = '  -1.48%  '
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.64'
$ws.Range("E13").Value = '  -1.12%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.875.25'
$ws.Range("E14").Value = '  -1.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.003'
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.179'
$ws.Range("E16").Value = '  -1.93%  '
$ws.Range("E17").Value = '  -2.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.83'
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06632'
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.04'
$ws.Range("E20").Value = '  +0.43%  '
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.075'
$ws.Range("E22").Value = '  -2.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.071.35'
$ws.Range("E23").Value = '  -0.57%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.41'
$ws.Range("E24").Value = '  -0.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.265'
$ws.Range("E25").Value = '  -2.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.577'
$ws.Range("E26").Value = '  -0.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.095.03'
$ws.Range("E27").Value = '  -1.39%  '
$ws.Range("E28").Value = '  -1.89%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '156.97'
$ws.Range("E29").Value = '  -0.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.13'
$ws.Range("E30").Value = '  -0.79%  '
$ws.Range("E31").Value = '  -2.31%  '
$ws.Range("E32").Value = '  -4.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.595'
$ws.Range("E33").Value = '  -0.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.598'
$ws.Range("E34").Value = '  -0.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.590'
$ws.Range("E35").Value = '  -1.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02439'
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06565'
$ws.Range("E37").Value = '  -1.56%  '
$ws.Range("E38").Value = '  -1.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.206'
$ws.Range("E39").Value = '  -2.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.239'
$ws.Range("E40").Value = '  -3.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.51'
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6368'
$ws.Range("E42").Value = '  -1.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.879'
$ws.Range("E43").Value = '  -2.43%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.19'
$ws.Range("E45").Value = '  -1.32%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5989'
$ws.Range("E46").Value = '  -1.56%  '
$ws.Range("B47").Value = 'WEMIXTOKEN'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.279'
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.674'
$ws.Range("E48").Value = '  -1.26%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.226'
$ws.Range("E49").Value = '  +3.25%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.986'
$ws.Range("E50").Value = '  -1.48%  '
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.20'
$ws.Range("E51").Value = '  -1.90%  '
